$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row after the last used row in column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "2025/12/05 02:00"
$ws.Cells.Item($newRow, 2).Value = "-"
$ws.Cells.Item($newRow, 3).Value = "-"
$ws.Cells.Item($newRow, 4).Value = "-"
$ws.Cells.Item($newRow, 5).Value = "-"
$ws.Cells.Item($newRow, 6).Value = "-"
$ws.Cells.Item($newRow, 7).Value = "-"
